$wb = $excel.ActiveWorkbook

# --- "tags" sheet: add a new tag row for the BioSharing reference ---
$tags = $wb.Worksheets.Item("tags")

$tags.Range("A4").Value = "miameenv_biosharing"
$tags.Range("B4").Value = "http://www.biosharing.org/bsg-000168"
$tags.Range("C4").Value = "http://www.biosharing.org/bsg-000168"
$tags.Range("D4").Value = "Reference"
$tags.Range("E4").Value = "system"
$tags.Range("F4").Value = "http://molgenis.org/biobankconnect/link"

# Match formatting of the row above (style index used by the rest of the table)
$tags.Range("A3:F3").Copy()
$tags.Range("A4:F4").PasteSpecial(-4122)

# Column widths (best-fit sizing for the new, wider URL content)
$tags.Columns.Item(1).ColumnWidth = 18.5
$tags.Range("B1:C1").EntireColumn.ColumnWidth = 46.833333333333336

# Landscape page setup (now needed once the sheet has a 6th, wider column of data)
$tags.PageSetup.Orientation = 2
$tags.PageSetup.PaperSize = 9

# --- "packages" sheet: reference the new tag from the MIAME-ENV package ---
$packages = $wb.Worksheets.Item("packages")
$packages.Range("D2").Value = "miameenv_home,miameenv_pub1,miameenv_biosharing"

# --- Active sheet / selection moves to "tags", cell D4 ---
$tags.Activate()
$tags.Range("D4").Select()
